$d = $word.ActiveDocument

$replacements = @(
    @("11×14=", "45×32="),
    @("87×79=", "28×81="),
    @("57×24=", "77×64="),
    @("59×71=", "44×49="),
    @("75×77=", "21×22="),
    @("72×33=", "14×76="),
    @("57×48=", "63×73="),
    @("78×76=", "45×87="),
    @("57×53=", "85×90="),
    @("15×41=", "38×54="),
    @("72×63=", "43×83="),
    @("37×31=", "28×38="),
    @("90×98=", "14×33="),
    @("33×92=", "81×40="),
    @("45×46=", "55×41="),
    @("47×67=", "49×20="),
    @("54×45=", "19×15="),
    @("51×91=", "56×22="),
    @("57×11=", "77×71="),
    @("26×75=", "30×72="),
    @("65×65=", "23×36="),
    @("49×90=", "71×56="),
    @("50×14=", "84×15="),
    @("38×30=", "73×60="),
    @("55×76=", "72×21=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
